# Update the cryptos price/volume table (rows 2-51) with the latest
# scraped values. Columns: B=Coin, C=Link, D=Price, E=Volume(1h).
# Numeric-looking price strings get their cell pre-formatted as Text
# ("@") before the write so Excel stores them verbatim (e.g. "0.999")
# instead of silently coercing them to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.921.20'
$ws.Range("E2").Value = '  +4.10%  '
$ws.Range("D3").Value = '2.348.14'
$ws.Range("E3").Value = '  +3.28%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '519.38'
$ws.Range("E5").Value = '  +3.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.46'
$ws.Range("E6").Value = '  +4.90%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.537'
$ws.Range("E8").Value = '  +1.98%  '
$ws.Range("D9").Value = '2.345.41'
$ws.Range("E9").Value = '  +2.53%  '
$ws.Range("E10").Value = '  +7.83%  '
$ws.Range("E11").Value = '  -0.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.24'
$ws.Range("E12").Value = '  +7.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.342'
$ws.Range("E13").Value = '  +0.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.76'
$ws.Range("E14").Value = '  +2.32%  '
$ws.Range("D15").Value = '2.739.64'
$ws.Range("E15").Value = '  +2.36%  '
$ws.Range("D16").Value = '56.820.94'
$ws.Range("E16").Value = '  +3.96%  '
$ws.Range("E17").Value = '  +3.05%  '
$ws.Range("D18").Value = '2.352.13'
$ws.Range("E18").Value = '  +2.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.45'
$ws.Range("E19").Value = '  +1.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.26'
$ws.Range("E20").Value = '  +3.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.03'
$ws.Range("E21").Value = '  +5.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.57'
$ws.Range("E22").Value = '  +2.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.73'
$ws.Range("E24").Value = '  +0.81%  '
$ws.Range("B25").Value = 'Kaspa'
$ws.Range("C25").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.161'
$ws.Range("E25").Value = '  +7.75%  '
$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.997'
$ws.Range("E26").Value = '  +0.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.77'
$ws.Range("E27").Value = '  +4.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.23'
$ws.Range("E28").Value = '  +11.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.39'
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("D30").Value = '0.0₃0743'
$ws.Range("E30").Value = '  +6.45%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.69'
$ws.Range("E31").Value = '  +4.10%  '
$ws.Range("E32").Value = '  +2.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.27'
$ws.Range("E33").Value = '  +2.08%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.25'
$ws.Range("E36").Value = '  +4.59%  '
$ws.Range("B37").Value = 'SuiNetwork'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.935'
$ws.Range("E37").Value = '  +3.21%  '
$ws.Range("E38").Value = '  +6.01%  '
$ws.Range("E39").Value = '  +8.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.71'
$ws.Range("E40").Value = '  +3.52%  '
$ws.Range("E41").Value = '  +1.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.60'
$ws.Range("E42").Value = '  +6.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '137.45'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '275.62'
$ws.Range("E44").Value = '  +10.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.09'
$ws.Range("E45").Value = '  +1.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0933'
$ws.Range("E46").Value = '  +3.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0505'
$ws.Range("E47").Value = '  +1.96%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.560'
$ws.Range("E48").Value = '  +2.63%  '
$ws.Range("E49").Value = '  +5.65%  '
$ws.Range("E50").Value = '  +1.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.86'
$ws.Range("E51").Value = '  +2.98%  '
